$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 2806.7693
$ws.Range("I34").Value = 2806.7693
$ws.Range("K34").Value = 2806.7693
$ws.Range("M34").Value = -2603.7693
$ws.Range("H36").Value = 2806.7693
$ws.Range("I36").Value = 2806.7693
$ws.Range("K36").Value = 2806.7693
$ws.Range("M36").Value = -2091.7693
$ws.Range("H51").Value = 9457.625
$ws.Range("I51").Value = 8489
$ws.Range("J51").Value = 9522.200000000001
$ws.Range("K51").Value = 8489
$ws.Range("L51").Value = 9522.200000000001
$ws.Range("M51").Value = -8005
$ws.Range("N51").Value = -10490.2
$ws.Range("H62").Value = 3054.0667
$ws.Range("I62").Value = 2321.9
$ws.Range("J62").Value = 4518.4
$ws.Range("K62").Value = 2321.9
$ws.Range("L62").Value = 4518.4
$ws.Range("M62").Value = -1697.9
$ws.Range("N62").Value = -5766.4
$ws.Range("H65").Value = 3054.0667
$ws.Range("I65").Value = 2321.9
$ws.Range("J65").Value = 4518.4
$ws.Range("K65").Value = 11609.5
$ws.Range("L65").Value = 22592
$ws.Range("M65").Value = -8489.5
$ws.Range("N65").Value = -28832
$ws.Range("H116").Value = 274518.7
$ws.Range("I116").Value = 557498.6
$ws.Range("J116").Value = 6432.4736
$ws.Range("K116").Value = 557498.6
$ws.Range("L116").Value = 6432.4736
$ws.Range("M116").Value = -554056.6
$ws.Range("N116").Value = -13316.4736
$ws.Range("H137").Value = 3666001
$ws.Range("I137").Value = 7938343
$ws.Range("K137").Value = 23815029
$ws.Range("M137").Value = -23812479
$ws.Range("H141").Value = 22519.2
$ws.Range("I141").Value = 26899
$ws.Range("K141").Value = 80697
$ws.Range("M141").Value = -75517

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 4308.222
$ws.Range("I36").Value = 396.2857
$ws.Range("J36").Value = 18000
$ws.Range("K36").Value = 396.2857
$ws.Range("L36").Value = 18000
$ws.Range("M36").Value = 137.7143
$ws.Range("N36").Value = -19068
$ws.Range("H115").Value = 31710.525
$ws.Range("J115").Value = 31710.525
$ws.Range("L115").Value = 31710.525
$ws.Range("N115").Value = -34844.525
$ws.Range("H122").Value = 41783.332
$ws.Range("J122").Value = 41783.332
$ws.Range("L122").Value = 41783.332
$ws.Range("N122").Value = -51583.332
$ws.Range("H134").Value = 2615.5
$ws.Range("I134").Value = 1928.8
$ws.Range("K134").Value = 5786.4
$ws.Range("M134").Value = -3251.4

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2721.742
$ws.Range("I31").Value = 1017.1905
$ws.Range("J31").Value = 6301.3
$ws.Range("K31").Value = 1017.1905
$ws.Range("L31").Value = 6301.3
$ws.Range("M31").Value = -722.1905
$ws.Range("N31").Value = -6891.3
$ws.Range("H34").Value = 2721.742
$ws.Range("I34").Value = 1017.1905
$ws.Range("J34").Value = 6301.3
$ws.Range("K34").Value = 1017.1905
$ws.Range("L34").Value = 6301.3
$ws.Range("M34").Value = -815.1905
$ws.Range("N34").Value = -6705.3
$ws.Range("H35").Value = 14432.9375
$ws.Range("I35").Value = 1393.3
$ws.Range("J35").Value = 36165.668
$ws.Range("K35").Value = 1393.3
$ws.Range("L35").Value = 36165.668
$ws.Range("M35").Value = -1099.3
$ws.Range("N35").Value = -36753.668
$ws.Range("H58").Value = 2776.1287
$ws.Range("I58").Value = 1675.1455
$ws.Range("J58").Value = 6813.067
$ws.Range("K58").Value = 1675.1455
$ws.Range("L58").Value = 6813.067
$ws.Range("M58").Value = -1472.1455
$ws.Range("N58").Value = -7219.067
$ws.Range("H136").Value = 2776.1287
$ws.Range("I136").Value = 1675.1455
$ws.Range("J136").Value = 6813.067
$ws.Range("K136").Value = 5025.4365
$ws.Range("L136").Value = 20439.201
$ws.Range("M136").Value = -2475.4365
$ws.Range("N136").Value = -25539.201

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 1100
$ws.Range("I36").Value = 200
$ws.Range("J36").Value = 2000
$ws.Range("K36").Value = 600
$ws.Range("L36").Value = 6000
$ws.Range("M36").Value = -431
$ws.Range("N36").Value = -6338

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 9259918
$ws.Range("I107").Value = 232.8
$ws.Range("J107").Value = 15873979
$ws.Range("K107").Value = 232.8
$ws.Range("L107").Value = 15873979
$ws.Range("M107").Value = 1687.2
$ws.Range("N107").Value = -15877819
$ws.Range("H126").Value = 4007.75
$ws.Range("I126").Value = 1745.5555
$ws.Range("K126").Value = 5236.666499999999
$ws.Range("M126").Value = -2766.666499999999
$ws.Range("H132").Value = 6033.231
$ws.Range("I132").Value = 5133.4287
$ws.Range("K132").Value = 15400.2861
$ws.Range("M132").Value = -12870.2861
$ws.Range("H134").Value = 38853.617
$ws.Range("J134").Value = 38853.617
$ws.Range("L134").Value = 116560.851
$ws.Range("N134").Value = -121630.851

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1852.7778
$ws.Range("J22").Value = 2428.5715
$ws.Range("L22").Value = 2428.5715
$ws.Range("N22").Value = -3018.5715
$ws.Range("H27").Value = 1852.7778
$ws.Range("J27").Value = 2428.5715
$ws.Range("L27").Value = 2428.5715
$ws.Range("N27").Value = -2642.5715
$ws.Range("H31").Value = 9980.909
$ws.Range("I31").Value = 1223.75
$ws.Range("J31").Value = 33333.332
$ws.Range("K31").Value = 1223.75
$ws.Range("L31").Value = 33333.332
$ws.Range("M31").Value = -975.75
$ws.Range("N31").Value = -33829.332
$ws.Range("H68").Value = 695.4375
$ws.Range("I68").Value = 695.4375
$ws.Range("K68").Value = 695.4375
$ws.Range("M68").Value = 53.5625
$ws.Range("H71").Value = 695.4375
$ws.Range("I71").Value = 695.4375
$ws.Range("K71").Value = 3477.1875
$ws.Range("M71").Value = 266.8125

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 15000
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 15000
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 15000
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -15470
$ws.Range("H26").Value = 7615
$ws.Range("I26").Value = 5012
$ws.Range("J26").Value = 8916.5
$ws.Range("K26").Value = 5012
$ws.Range("L26").Value = 8916.5
$ws.Range("M26").Value = -4719
$ws.Range("N26").Value = -9502.5
$ws.Range("H35").Value = 15000
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 15000
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 15000
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -30530
